$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: 'Bitcoin' -> 'Bitcoin'
$ws.Range("D2").Value = "'36.660.72"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.72%  "
$ws.Range("E2").Style = "Normal"

# Row 3: 'Ethereum' -> 'Ethereum'
$ws.Range("D3").Value = "'2.050.80"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.61%  "
$ws.Range("E3").Style = "Normal"

# Row 4: 'TetherUSD' -> 'TetherUSD'
$ws.Range("E4").Value = "'  -0.05%  "
$ws.Range("E4").Style = "Normal"

# Row 5: 'BNB' -> 'BNB'
$ws.Range("D5").Value = "'245.59"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.15%  "
$ws.Range("E5").Style = "Normal"

# Row 6: 'XRP' -> 'XRP'
$ws.Range("D6").Value = "'0.669"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +2.15%  "
$ws.Range("E6").Style = "Normal"

# Row 7: 'USDC' -> 'Solana'
$ws.Range("B7").Value = "Solana"
$ws.Range("C7").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D7").Value = "'57.79"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -2.38%  "
$ws.Range("E7").Style = "Normal"

# Row 8: 'Solana' -> 'USDC'
$ws.Range("B8").Value = "USDC"
$ws.Range("C8").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D8").Value = "'1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +0.02%  "
$ws.Range("E8").Style = "Normal"

# Row 9: 'OKB' -> 'OKB'
$ws.Range("D9").Value = "'63.28"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +6.92%  "
$ws.Range("E9").Style = "Normal"

# Row 10: 'Cardano' -> 'Cardano'
$ws.Range("D10").Value = "'0.369"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -0.84%  "
$ws.Range("E10").Style = "Normal"

# Row 11: 'Dogecoin' -> 'Dogecoin'
$ws.Range("D11").Value = "'0.0751"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -3.91%  "
$ws.Range("E11").Style = "Normal"

# Row 12: 'TRON' -> 'TRON'
$ws.Range("D12").Value = "'0.107"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -3.32%  "
$ws.Range("E12").Style = "Normal"

# Row 13: 'Polygon' -> 'Polygon'
$ws.Range("E13").Value = "'  +3.58%  "
$ws.Range("E13").Style = "Normal"

# Row 14: 'Chainlink' -> 'Chainlink'
$ws.Range("D14").Value = "'14.73"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -3.23%  "
$ws.Range("E14").Style = "Normal"

# Row 15: 'WrappedliquidstakedEther2.0' -> 'WrappedliquidstakedEther2.0'
$ws.Range("D15").Value = "'2.346.58"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -0.67%  "
$ws.Range("E15").Style = "Normal"

# Row 16: 'Polkadot' -> 'Polkadot'
$ws.Range("D16").Value = "'5.44"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -2.92%  "
$ws.Range("E16").Style = "Normal"

# Row 17: 'WrappedEther' -> 'WrappedEther'
$ws.Range("D17").Value = "'2.035.34"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -2.25%  "
$ws.Range("E17").Style = "Normal"

# Row 18: 'WrappedBTC' -> 'Avalanche'
$ws.Range("B18").Value = "Avalanche"
$ws.Range("C18").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D18").Value = "'17.62"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +0.34%  "
$ws.Range("E18").Style = "Normal"

# Row 19: 'Avalanche' -> 'WrappedBTC'
$ws.Range("B19").Value = "WrappedBTC"
$ws.Range("C19").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D19").Value = "'36.524.99"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -1.03%  "
$ws.Range("E19").Style = "Normal"

# Row 20: 'Litecoin' -> 'Litecoin'
$ws.Range("D20").Value = "'72.19"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -1.83%  "
$ws.Range("E20").Style = "Normal"

# Row 21: 'ShibaInu' -> 'ShibaInu'
$ws.Range("D21").Value = "'0.0₃0860"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -3.73%  "
$ws.Range("E21").Style = "Normal"

# Row 22: 'BitcoinCash' -> 'BitcoinCash'
$ws.Range("D22").Value = "'238.47"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +0.86%  "
$ws.Range("E22").Style = "Normal"

# Row 23: 'Uniswap' -> 'Uniswap'
$ws.Range("D23").Value = "'5.20"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -4.59%  "
$ws.Range("E23").Style = "Normal"

# Row 24: 'Dai' -> 'Dai'
$ws.Range("E24").Value = "'  +0.04%  "
$ws.Range("E24").Style = "Normal"

# Row 25: 'Toncoin' -> 'Toncoin'
$ws.Range("E25").Value = "'  -2.60%  "
$ws.Range("E25").Style = "Normal"

# Row 26: 'PancakeSwap' -> 'PancakeSwap'
$ws.Range("D26").Value = "'2.28"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +3.07%  "
$ws.Range("E26").Style = "Normal"

# Row 27: 'Cosmos' -> 'Cosmos'
$ws.Range("D27").Value = "'9.32"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -6.98%  "
$ws.Range("E27").Style = "Normal"

# Row 28: 'Monero' -> 'Monero'
$ws.Range("D28").Value = "'164.60"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -2.47%  "
$ws.Range("E28").Style = "Normal"

# Row 29: 'EthereumClassic' -> 'EthereumClassic'
$ws.Range("D29").Value = "'20.09"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -0.93%  "
$ws.Range("E29").Style = "Normal"

# Row 30: 'Stellar' -> 'Stellar'
$ws.Range("D30").Value = "'0.122"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -1.91%  "
$ws.Range("E30").Style = "Normal"

# Row 31: 'ImmutableX' -> 'ImmutableX'
$ws.Range("D31").Value = "'1.19"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +3.67%  "
$ws.Range("E31").Style = "Normal"

# Row 32: 'Filecoin' -> 'Filecoin'
$ws.Range("D32").Value = "'5.04"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -7.92%  "
$ws.Range("E32").Style = "Normal"

# Row 33: 'Hedera' -> 'Hedera'
$ws.Range("D33").Value = "'0.0602"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -2.24%  "
$ws.Range("E33").Style = "Normal"

# Row 34: 'InternetComputer(DFINITY)' -> 'InternetComputer(DFINITY)'
$ws.Range("D34").Value = "'4.44"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -8.30%  "
$ws.Range("E34").Style = "Normal"

# Row 35: 'Kaspa' -> 'Kaspa'
$ws.Range("D35").Value = "'0.0869"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +1.74%  "
$ws.Range("E35").Style = "Normal"

# Row 36: 'BinanceUSD' -> 'BinanceUSD'
$ws.Range("E36").Value = "'  -0.09%  "
$ws.Range("E36").Style = "Normal"

# Row 37: 'WEMIXToken' -> 'WEMIXToken'
$ws.Range("E37").Value = "'  -1.29%  "
$ws.Range("E37").Style = "Normal"

# Row 38: 'LidoDAOToken' -> 'LidoDAOToken'
$ws.Range("D38").Value = "'2.21"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -6.30%  "
$ws.Range("E38").Style = "Normal"

# Row 39: 'THORChain' -> 'THORChain'
$ws.Range("D39").Value = "'5.03"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +2.08%  "
$ws.Range("E39").Style = "Normal"

# Row 40: 'TrustWalletToken' -> 'TrustWalletToken'
$ws.Range("E40").Value = "'  -6.09%  "
$ws.Range("E40").Style = "Normal"

# Row 41: 'VeChain' -> 'HuobiToken'
$ws.Range("B41").Value = "HuobiToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D41").Value = "'2.92"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -1.37%  "
$ws.Range("E41").Style = "Normal"

# Row 42: 'HuobiToken' -> 'VeChain'
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").Value = "'0.0216"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -2.89%  "
$ws.Range("E42").Style = "Normal"

# Row 43: 'ARBITRUM' -> 'ARBITRUM'
$ws.Range("D43").Value = "'1.11"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -4.05%  "
$ws.Range("E43").Style = "Normal"

# Row 44: 'Aave' -> 'Aave'
$ws.Range("D44").Value = "'94.28"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -2.98%  "
$ws.Range("E44").Style = "Normal"

# Row 45: 'Cronos' -> 'Cronos'
$ws.Range("D45").Value = "'0.0907"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -6.54%  "
$ws.Range("E45").Style = "Normal"

# Row 46: 'InjectiveProtocol' -> 'InjectiveProtocol'
$ws.Range("D46").Value = "'16.02"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -4.51%  "
$ws.Range("E46").Style = "Normal"

# Row 47: 'FraxShare' -> 'Maker'
$ws.Range("B47").Value = "Maker"
$ws.Range("C47").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D47").Value = "'1.381.60"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +4.77%  "
$ws.Range("E47").Style = "Normal"

# Row 48: 'Maker' -> 'FraxShare'
$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D48").Value = "'7.48"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +8.90%  "
$ws.Range("E48").Style = "Normal"

# Row 49: 'MXToken' -> 'MXToken'
$ws.Range("D49").Value = "'2.94"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +2.39%  "
$ws.Range("E49").Style = "Normal"

# Row 50: 'RenderToken' -> 'RenderToken'
$ws.Range("D50").Value = "'2.26"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -5.09%  "
$ws.Range("E50").Style = "Normal"

# Row 51: 'MultiversX' -> 'MultiversX'
$ws.Range("D51").Value = "'45.79"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +0.86%  "
$ws.Range("E51").Style = "Normal"

